## NIT-9003429496.xlsx — "Actualiza base de datos EC y agrega parte 1 de
## nuevos estado de cuenta"
##
## The worker/debtor table (rows 16-34) is replaced with a new, shorter
## data set (rows 16-22), and the signature footer (previously rows 39-40)
## shifts up to sit right below it (rows 27-28). The account-level summary
## values (E11, C13, F13) are also updated.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# ---------------------------------------------------------------------
# 1) Remove the 12 rows of data that are no longer needed. Deleting the
#    FIRST 12 data rows (16-27) -- rather than the last 12 -- keeps the
#    special "bottom border" row style (currently on row 34) attached to
#    its row, which then slides up to become the new last data row (22).
# ---------------------------------------------------------------------
$ws.Range("A16:A27").EntireRow.Delete() | Out-Null

# ---------------------------------------------------------------------
# 2) Update the account summary block.
# ---------------------------------------------------------------------
$ws.Range("E11").Value = 176733
$ws.Range("C13").Value = 5
$ws.Range("F13").Value = 3

# ---------------------------------------------------------------------
# 3) Overwrite the (now 7-row) worker/debtor table with the new data.
#    Columns: B=Tipo Doc, C=N Doc, D=Nombre, E=Periodo Mora,
#             F=Valor Mora, G=Salario Basico
# ---------------------------------------------------------------------
$rows = @(
    @("CC", "9156419",    "RAFAEL CIPRIANO VALDES OSORIO",   "1912", 24292, 877803),
    @("CC", "9043421",    "EDINSON JACOME CONTRERAS",        "1912",  2208, 877803),
    @("CC", "1104413446", "DUVAN MANUEL MONTERROZA JARABA",  "1912", 24292, 828116),
    @("CC", "9156419",    "RAFAEL CIPRIANO VALDES OSORIO",   "2001", 33125, 877803),
    @("CC", "9156419",    "RAFAEL CIPRIANO VALDES OSORIO",   "2002", 33125, 877803),
    @("CC", "10965261",   "JOSE LUIS SERPA MEJIA",           "2002", 28090, 877803),
    @("CC", "91279183",   "JUAN CARLOS MONSALVE MARTINEZ",   "2002", 31601, 877803)
)

$r = 16
foreach ($row in $rows) {
    $ws.Range("B$r").Value = $row[0]
    $ws.Range("C$r").Value = $row[1]
    $ws.Range("D$r").Value = $row[2]
    $ws.Range("E$r").Value = $row[3]
    $ws.Range("F$r").Value = $row[4]
    $ws.Range("G$r").Value = $row[5]
    $r++
}

# ---------------------------------------------------------------------
# 4) Column D ("Nombre Trabajador") is best-fit width; with the longest
#    name now shorter, re-fit it to the new content.
# ---------------------------------------------------------------------
$ws.Columns.Item(4).AutoFit() | Out-Null

Write-Host "done"
